# Update NATMI LR-pair TPM-derived values per new computation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Psap / ECs -> MuSCs)
$ws.Range("G2").Value = 137.0025913333334
$ws.Range("H2").Value = 411.007774
$ws.Range("I2").Value = 0.07043159922291199
$ws.Range("J2").Value = 0.07043159922291199
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.06874933333333333
$ws.Range("N2").Value = 0.206248
$ws.Range("Q2").Value = 9.418836819105779
$ws.Range("R2").Value = 84.769531371952
$ws.Range("S2").Value = 0.07043159922291199
$ws.Range("T2").Value = 0.07043159922291199

# Row 3 (Psap / FAPs -> MuSCs)
$ws.Range("I3").Value = 0.1284841594777439
$ws.Range("J3").Value = 0.1284841594777439
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.06874933333333333
$ws.Range("N3").Value = 0.206248
$ws.Range("Q3").Value = 17.182221407904
$ws.Range("R3").Value = 154.639992671136
$ws.Range("S3").Value = 0.1284841594777439
$ws.Range("T3").Value = 0.1284841594777439

# Row 4 (Psap / MuSCs -> MuSCs)
$ws.Range("G4").Value = 88.73577866666666
$ws.Range("H4").Value = 266.207336
$ws.Range("I4").Value = 0.04561813567874526
$ws.Range("J4").Value = 0.04561813567874527
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.06874933333333333
$ws.Range("N4").Value = 0.206248
$ws.Range("Q4").Value = 6.100525626147554
$ws.Range("R4").Value = 54.90473063532799
$ws.Range("S4").Value = 0.04561813567874526
$ws.Range("T4").Value = 0.04561813567874527

# Row 5 (Psap / Resolving-Mac -> MuSCs)
$ws.Range("G5").Value = 1469.52242
$ws.Range("H5").Value = 4408.56726
$ws.Range("I5").Value = 0.7554661056205989
$ws.Range("J5").Value = 0.7554661056205988
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.06874933333333333
$ws.Range("N5").Value = 0.206248
$ws.Range("Q5").Value = 101.0286866933867
$ws.Range("R5").Value = 909.2581802404799
$ws.Range("S5").Value = 0.7554661056205989
$ws.Range("T5").Value = 0.7554661056205988
